# Bugfixed the naive forecaster component module
#
# The underlying data source dropped its earliest observation (the row that
# used to be row 2), so every row shifts up by one. The y_1_forecast column
# (E) is recomputed from scratch against the corrected series, and the C3
# value is refreshed with its newly computed (slightly different) value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the obsolete first data row; Excel shifts rows 3:19 up to 2:18
# and keeps everything else (styles, formatting, columns A/B/C/D values,
# dimension) correctly adjusted automatically.
$ws.Rows("2:2").Delete()

# Refresh the recomputed y_0_forecast value for 2009 (row 3).
$ws.Range("C3").Value = 0.1715429114845124

# Refresh the recomputed y_1_forecast column (E) for every remaining row.
$eValues = @{
    2  = $null
    3  = $null
    4  = $null
    5  = $null
    6  = 0.5784444854042281
    7  = 2.529895848567842
    8  = 4.060884847379076
    9  = 2.270469368501771
    10 = 2.467161166346266
    11 = 2.480855794925163
    12 = 3.221757900820066
    13 = 2.631992339577627
    14 = 2.153309886824961
    15 = 4.667362054855917
    16 = 3.305715257492858
    17 = 1.757655717321982
    18 = 2.159361127638926
}

foreach ($row in $eValues.Keys) {
    $val = $eValues[$row]
    $cell = $ws.Range("E$row")
    if ($null -eq $val) {
        $cell.ClearContents()
    } else {
        $cell.Value = $val
    }
}
